$d = $word.ActiveDocument

# Locate the paragraph that begins "Please read Using Pressure Canners..."
# and the empty separator paragraph ("\n") immediately preceding it, then
# delete both paragraphs outright. This removes the "blank line" + the
# "Please read Using Pressure Canners before beginning..." paragraph from
# between "Select young, tender pods..." and the next blank-line
# separator, per the target revision.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Please read Using Pressure Canners*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ge 2) {
    $target = $d.Paragraphs.Item($targetIndex)
    $prev = $d.Paragraphs.Item($targetIndex - 1)
    $prevTxt = $prev.Range.Text.TrimEnd([char]13)

    # Delete the "Please read..." paragraph first (higher in the document
    # order doesn't matter here since it's a single paragraph), then the
    # blank separator paragraph before it, so indices stay valid.
    $target.Range.Delete()
    if ($prevTxt -eq "\n") {
        $prev.Range.Delete()
    }
}

Write-Output ("Paragraphs remaining: " + $d.Paragraphs.Count)
